$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "As the system, I reject invalid transactions" (row 6) as DONE
$ws.Range("C6").Value = "DONE"

# Move the active selection to C12 (matches author's last-saved cursor position)
$ws.Range("C12").Select()
